$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "D2" = '28.387.82'
    "E2" = '  +1.36%  '
    "D3" = '1.879.19'
    "E3" = '  +0.83%  '
    "D4" = '1.011'
    "E4" = '  +0.74%  '
    "D5" = '315.63'
    "E5" = '  +1.13%  '
    "D6" = '1.012'
    "E6" = '  +0.95%  '
    "D7" = '0.5137'
    "E7" = '  +0.17%  '
    "D8" = '0.3935'
    "E8" = '  +1.22%  '
    "D9" = '0.08317'
    "E9" = '  -0.20%  '
    "D10" = '1.122'
    "B11" = 'Polkadot'
    "C11" = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
    "D11" = '6.280'
    "E11" = '  +2.09%  '
    "B12" = 'WrappedEther'
    "C12" = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
    "D12" = '1.877.58'
    "E12" = '  +0.68%  '
    "B13" = 'Solana'
    "C13" = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
    "D13" = '20.40'
    "E13" = '  -0.05%  '
    "B14" = 'Chainlink'
    "C14" = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
    "D14" = '7.253'
    "E14" = '  +0.04%  '
    "B15" = 'BinanceUSD'
    "C15" = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
    "D15" = '1.011'
    "E15" = '  +0.72%  '
    "B16" = 'ShibaInu'
    "C16" = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
    "D16" = '0.00001107'
    "E16" = '  +1.02%  '
    "B17" = 'Litecoin'
    "C17" = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
    "D17" = '91.32'
    "E17" = '  +0.83%  '
    "B18" = 'TRON'
    "C18" = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
    "D18" = '0.06738'
    "E18" = '  +1.84%  '
    "B19" = 'Avalanche'
    "C19" = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
    "D19" = '17.75'
    "E19" = '  +0.96%  '
    "B20" = 'Dai'
    "C20" = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
    "D20" = '1.013'
    "E20" = '  +1.08%  '
    "B21" = 'Uniswap'
    "C21" = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
    "D21" = '6.010'
    "E21" = '  +0.39%  '
    "B22" = 'WrappedBTC'
    "C22" = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
    "D22" = '28.432.55'
    "E22" = '  +1.38%  '
    "B23" = 'Cosmos'
    "C23" = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
    "D23" = '11.16'
    "E23" = '  +1.30%  '
    "B24" = 'Toncoin'
    "C24" = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
    "D24" = '2.247'
    "E24" = '  +0.22%  '
    "B25" = 'WrappedliquidstakedEther2.0'
    "C25" = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
    "D25" = '2.089.63'
    "E25" = '  +0.21%  '
    "B26" = 'Monero'
    "C26" = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
    "D26" = '160.75'
    "E26" = '  +1.69%  '
    "B27" = 'EthereumClassic'
    "C27" = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
    "D27" = '20.80'
    "E27" = '  +1.59%  '
    "B28" = 'LidoDAOToken'
    "C28" = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
    "D28" = '2.448'
    "E28" = '  -0.31%  '
    "B29" = 'BitcoinCash'
    "C29" = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
    "D29" = '127.11'
    "E29" = '  +1.61%  '
    "B30" = 'Stellar'
    "C30" = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
    "D30" = '0.1065'
    "E30" = '  +0.28%  '
    "B31" = 'ImmutableX'
    "C31" = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
    "D31" = '1.052'
    "E31" = '  +2.56%  '
    "B32" = 'Filecoin'
    "C32" = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
    "D32" = '5.898'
    "E32" = '  +1.14%  '
    "B33" = 'HuobiToken'
    "C33" = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
    "D33" = '3.639'
    "E33" = '  +1.31%  '
    "B34" = 'VeChain'
    "C34" = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
    "D34" = '0.02443'
    "E34" = '  +0.68%  '
    "B35" = 'Hedera'
    "C35" = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
    "D35" = '0.06532'
    "E35" = '  +0.18%  '
    "B36" = 'FraxShare'
    "C36" = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
    "D36" = '9.226'
    "E36" = '  -1.40%  '
    "B37" = 'Algorand'
    "C37" = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
    "D37" = '0.2191'
    "E37" = '  +0.79%  '
    "B38" = 'TrustWalletToken'
    "C38" = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
    "D38" = '1.258'
    "E38" = '  +3.62%  '
    "B39" = 'TheSandbox'
    "C39" = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
    "D39" = '0.6485'
    "E39" = '  +0.26%  '
    "D40" = '1.191'
    "E40" = '  -0.59%  '
    "B41" = 'InternetComputer(DFINITY)'
    "C41" = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
    "D41" = '4.992'
    "E41" = '  +0.66%  '
    "B42" = 'Aptos'
    "C42" = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
    "D42" = '11.20'
    "E42" = '  -0.89%  '
    "B43" = 'Decentraland'
    "C43" = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
    "D43" = '0.6070'
    "E43" = '  +0.19%  '
    "B44" = 'EnergySwap'
    "C44" = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
    "D44" = '13.17'
    "E44" = '  +1.60%  '
    "B45" = 'PancakeSwap'
    "C45" = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
    "D45" = '3.701'
    "E45" = '  +0.79%  '
    "B46" = 'WEMIXTOKEN'
    "C46" = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
    "D46" = '1.282'
    "E46" = '  -0.64%  '
    "B47" = 'NEARProtocol'
    "C47" = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
    "D47" = '2.025'
    "E47" = '  +1.57%  '
    "B48" = 'EOS'
    "C48" = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
    "D48" = '1.222'
    "E48" = '  +0.41%  '
    "B49" = 'Quant'
    "C49" = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
    "D49" = '122.12'
    "E49" = '  +1.07%  '
    "B50" = 'Cronos'
    "C50" = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
    "D50" = '0.06904'
    "E50" = '  +0.70%  '
    "B51" = 'Aave'
    "C51" = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
    "D51" = '77.78'
    "E51" = '  -0.57%  '
}

foreach ($addr in $values.Keys) {
    $cell = $ws.Range($addr)
    $cell.Value = "'" + $values[$addr]
    $cell.Style = "Normal"
}
